$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2087.1094
$ws.Range("I15").Value = 2087.1094
$ws.Range("K15").Value = 6261.3282
$ws.Range("M15").Value = -6092.3282
$ws.Range("H87").Value = 69999.336
$ws.Range("J87").Value = 69999.336
$ws.Range("L87").Value = 69999.336
$ws.Range("N87").Value = -72495.336
$ws.Range("H90").Value = 69999.336
$ws.Range("J90").Value = 69999.336
$ws.Range("L90").Value = 209998.008
$ws.Range("N90").Value = -222478.008
$ws.Range("H98").Value = 2380.15
$ws.Range("I98").Value = 1673.5333
$ws.Range("K98").Value = 1673.5333
$ws.Range("M98").Value = -175.5333000000001
$ws.Range("H101").Value = 354.75
$ws.Range("I101").Value = 354.75
$ws.Range("K101").Value = 1064.25
$ws.Range("M101").Value = 557.75
$ws.Range("H111").Value = 876
$ws.Range("I111").Value = 708.5714
$ws.Range("K111").Value = 2125.7142
$ws.Range("M111").Value = 941.2857999999997
$ws.Range("H113").Value = 6143.6665
$ws.Range("I113").Value = 6143.6665
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 6143.6665
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -2889.6665
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 2380.15
$ws.Range("I122").Value = 1673.5333
$ws.Range("K122").Value = 5020.5999
$ws.Range("M122").Value = -2570.5999
$ws.Range("H138").Value = 3025.0925
$ws.Range("J138").Value = 3398.9487
$ws.Range("L138").Value = 10196.8461
$ws.Range("N138").Value = -20476.8461

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2075.4167
$ws.Range("I2").Value = 2143.3635
$ws.Range("K2").Value = 2143.3635
$ws.Range("M2").Value = -2030.3635
$ws.Range("H32").Value = 188630.42
$ws.Range("I32").Value = 271445.62
$ws.Range("J32").Value = 23000
$ws.Range("K32").Value = 271445.62
$ws.Range("L32").Value = 23000
$ws.Range("M32").Value = -271158.62
$ws.Range("N32").Value = -23574
$ws.Range("H43").Value = 306223.25
$ws.Range("I43").Value = 74900
$ws.Range("J43").Value = 383331
$ws.Range("K43").Value = 74900
$ws.Range("L43").Value = 383331
$ws.Range("M43").Value = -74587
$ws.Range("N43").Value = -383957
$ws.Range("H45").Value = 1888.6154
$ws.Range("I45").Value = 1858.6666
$ws.Range("J45").Value = 1914.2858
$ws.Range("K45").Value = 1858.6666
$ws.Range("L45").Value = 1914.2858
$ws.Range("M45").Value = -1481.6666
$ws.Range("N45").Value = -2668.2858
$ws.Range("H61").Value = 1473573.4
$ws.Range("I61").Value = 2823.7666
$ws.Range("K61").Value = 2823.7666
$ws.Range("M61").Value = -2611.7666
$ws.Range("H74").Value = 1118621.2
$ws.Range("I74").Value = 1504960.9
$ws.Range("J74").Value = 19039.23
$ws.Range("K74").Value = 1504960.9
$ws.Range("L74").Value = 19039.23
$ws.Range("M74").Value = -1504086.9
$ws.Range("N74").Value = -20787.23
$ws.Range("H77").Value = 1118621.2
$ws.Range("I77").Value = 1504960.9
$ws.Range("J77").Value = 19039.23
$ws.Range("K77").Value = 7524804.5
$ws.Range("L77").Value = 95196.14999999999
$ws.Range("M77").Value = -7520436.5
$ws.Range("N77").Value = -103932.15
$ws.Range("H102").Value = 888.625
$ws.Range("I102").Value = 888.625
$ws.Range("K102").Value = 888.625
$ws.Range("M102").Value = 733.375
$ws.Range("H116").Value = 2075.4167
$ws.Range("I116").Value = 2143.3635
$ws.Range("K116").Value = 2143.3635
$ws.Range("M116").Value = 150.6365000000001
$ws.Range("H120").Value = 50000
$ws.Range("I120").Value = 50000
$ws.Range("K120").Value = 50000
$ws.Range("M120").Value = -45162
$ws.Range("H122").Value = 1869.1471
$ws.Range("I122").Value = 1829.4828
$ws.Range("K122").Value = 5488.4484
$ws.Range("M122").Value = -3038.4484
$ws.Range("H134").Value = 52784.9
$ws.Range("J134").Value = 52784.9
$ws.Range("L134").Value = 52784.9
$ws.Range("N134").Value = -62924.9
$ws.Range("H136").Value = 1473573.4
$ws.Range("I136").Value = 2823.7666
$ws.Range("K136").Value = 8471.299800000001
$ws.Range("M136").Value = -5921.299800000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2075.4167
$ws.Range("I3").Value = 2143.3635
$ws.Range("K3").Value = 2143.3635
$ws.Range("M3").Value = -2029.3635
$ws.Range("H35").Value = 39000
$ws.Range("J35").Value = 39000
$ws.Range("L35").Value = 39000
$ws.Range("N35").Value = -39620
$ws.Range("H37").Value = 1345.2
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H99").Value = 12073.643
$ws.Range("I99").Value = 16103.1
$ws.Range("K99").Value = 16103.1
$ws.Range("M99").Value = -14605.1
$ws.Range("H105").Value = 7543.032
$ws.Range("I105").Value = 6797.3335
$ws.Range("J105").Value = 10099.714
$ws.Range("K105").Value = 6797.3335
$ws.Range("L105").Value = 10099.714
$ws.Range("M105").Value = -5050.3335
$ws.Range("N105").Value = -13593.714
$ws.Range("H134").Value = 4514524
$ws.Range("I134").Value = 7596.357
$ws.Range("K134").Value = 22789.071
$ws.Range("M134").Value = -20254.071

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 168805.44
$ws.Range("I16").Value = 85208.164
$ws.Range("K16").Value = 85208.164
$ws.Range("M16").Value = -84921.164
$ws.Range("H113").Value = 168805.44
$ws.Range("I113").Value = 85208.164
$ws.Range("K113").Value = 85208.164
$ws.Range("M113").Value = -83038.164

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 1750
$ws.Range("I81").Value = 1000
$ws.Range("J81").Value = 2000
$ws.Range("K81").Value = 3000
$ws.Range("L81").Value = 6000
$ws.Range("M81").Value = -1877
$ws.Range("N81").Value = -8246
$ws.Range("H84").Value = 1750
$ws.Range("I84").Value = 1000
$ws.Range("J84").Value = 2000
$ws.Range("K84").Value = 9000
$ws.Range("L84").Value = 18000
$ws.Range("M84").Value = -3384
$ws.Range("N84").Value = -29232
$ws.Range("H109").Value = 6407.619
$ws.Range("I109").Value = 4163.5
$ws.Range("J109").Value = 8447.727999999999
$ws.Range("K109").Value = 12490.5
$ws.Range("L109").Value = 25343.184
$ws.Range("M109").Value = -11450.5
$ws.Range("N109").Value = -27423.184
$ws.Range("H132").Value = 8400.467000000001
$ws.Range("I132").Value = 1348.5
$ws.Range("K132").Value = 12136.5
$ws.Range("M132").Value = -9606.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 5097.385
$ws.Range("I21").Value = 5123.6665
$ws.Range("K21").Value = 5123.6665
$ws.Range("M21").Value = -4950.6665
$ws.Range("H30").Value = 5097.385
$ws.Range("I30").Value = 5123.6665
$ws.Range("K30").Value = 5123.6665
$ws.Range("M30").Value = -5018.6665
$ws.Range("H93").Value = 39874.5
$ws.Range("J93").Value = 39874.5
$ws.Range("L93").Value = 39874.5
$ws.Range("N93").Value = -43618.5
$ws.Range("H95").Value = 48672
$ws.Range("J95").Value = 48672
$ws.Range("L95").Value = 48672
$ws.Range("N95").Value = -54164
$ws.Range("H132").Value = 11667.521
$ws.Range("I132").Value = 10259.395
$ws.Range("J132").Value = 18356.125
$ws.Range("K132").Value = 30778.185
$ws.Range("L132").Value = 55068.375
$ws.Range("M132").Value = -28248.185
$ws.Range("N132").Value = -60128.375

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3449.3
$ws.Range("I22").Value = 2250
$ws.Range("J22").Value = 3749.125
$ws.Range("K22").Value = 2250
$ws.Range("L22").Value = 3749.125
$ws.Range("M22").Value = -1955
$ws.Range("N22").Value = -4339.125
$ws.Range("H27").Value = 3449.3
$ws.Range("I27").Value = 2250
$ws.Range("J27").Value = 3749.125
$ws.Range("K27").Value = 2250
$ws.Range("L27").Value = 3749.125
$ws.Range("M27").Value = -2143
$ws.Range("N27").Value = -3963.125
$ws.Range("H46").Value = 5249.3125
$ws.Range("I46").Value = 2001
$ws.Range("K46").Value = 2001
$ws.Range("M46").Value = -1813
$ws.Range("H93").Value = 1468.5807
$ws.Range("J93").Value = 1759.8
$ws.Range("L93").Value = 1759.8
$ws.Range("N93").Value = -4255.8
$ws.Range("H116").Value = 223822.67
$ws.Range("J116").Value = 223822.67
$ws.Range("L116").Value = 223822.67
$ws.Range("N116").Value = -233000.67
$ws.Range("H140").Value = 53666
$ws.Range("J140").Value = 57999
$ws.Range("L140").Value = 57999
$ws.Range("N140").Value = -68359

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 15950.956
$ws.Range("J62").Value = 15996.421
$ws.Range("L62").Value = 15996.421
$ws.Range("N62").Value = -17244.421
$ws.Range("H65").Value = 15950.956
$ws.Range("J65").Value = 15996.421
$ws.Range("L65").Value = 79982.105
$ws.Range("N65").Value = -86222.105
$ws.Range("H132").Value = 2925405.5
$ws.Range("I132").Value = 3334638.2
$ws.Range("K132").Value = 10003914.6
$ws.Range("M132").Value = -10001384.6
$ws.Range("H136").Value = 10587000
$ws.Range("I136").Value = 2289158
$ws.Range("K136").Value = 6867474
$ws.Range("M136").Value = -6864924
